$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 574, shifting existing rows 574:649 down to 575:650
$ws.Rows.Item(574).Insert()

# Populate the newly inserted row 574 with the new record's data
$ws.Range("A574").Value2 = 3
$ws.Range("B574").Value2 = "Femacal de La Calera"
$ws.Range("C574").Value2 = "Coquimbo"
$ws.Range("D574").Value2 = 44984
$ws.Range("E574").Value2 = 5
$ws.Range("F574").Value2 = 100112037
$ws.Range("G574").Value2 = "Cebollín"
$ws.Range("H574").Value2 = "Sin especificar"
$ws.Range("I574").Value2 = "Primera"
$ws.Range("J574").Value2 = 210
$ws.Range("K574").Value2 = 3500
$ws.Range("L574").Value2 = 4000
$ws.Range("M574").Value2 = 3738
$ws.Range("N574").Value2 = "$/paquete 36 unidades"
$ws.Range("O574").Value2 = "Provincia de Quillota"
$ws.Range("P574").Value2 = 104
$ws.Range("Q574").Value2 = 36
$ws.Range("R574").Value2 = "Hortaliza"

# Match the date-format style used by the rest of column D
$ws.Range("D574").NumberFormat = $ws.Range("D575").NumberFormat
